$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.008633180391209853
$ws.Range("C2").Value = 0.009875948452366614
$ws.Range("D2").Value = 0.009254564421788233
$ws.Range("E2").Value = 0.0006213840305783801

$ws.Range("B3").Value = 0.2395309882747069
$ws.Range("C3").Value = 0.2710743801652892
$ws.Range("D3").Value = 0.255302684219998
$ws.Range("E3").Value = 0.01577169594529118

$ws.Range("B4").Value = 0.01666569547229183
$ws.Range("C4").Value = 0.01905757945499971
$ws.Range("D4").Value = 0.01786163746364577
$ws.Range("E4").Value = 0.001195941991353939
